$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 43 (shifts existing rows 43-130 down to 45-132)
$ws.Rows("43:44").Insert()

# Row 43
$ws.Range("A43").Value = 9
$ws.Range('B43').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C43').Value = 'Metropolitana'
$ws.Range("D43").Value = 44469
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = 100112017
$ws.Range('G43').Value = 'Apio'
$ws.Range('H43').Value = 'Americana (o)'
$ws.Range('I43').Value = 'Primera'
$ws.Range("J43").Value = 61
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 9000
$ws.Range("M43").Value = 8492
$ws.Range('N43').Value = '$/docena de matas'
$ws.Range('O43').Value = 'Región de Coquimbo'
$ws.Range("P43").Value = 1415
$ws.Range("Q43").Value = 6
$ws.Range('R43').Value = 'Hortaliza'

# Row 44
$ws.Range("A44").Value = 9
$ws.Range('B44').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C44').Value = 'Metropolitana'
$ws.Range("D44").Value = 44469
$ws.Range("E44").Value = 13
$ws.Range("F44").Value = 100112017
$ws.Range('G44').Value = 'Apio'
$ws.Range('H44').Value = 'Americana (o)'
$ws.Range('I44').Value = 'Segunda'
$ws.Range("J44").Value = 34
$ws.Range("K44").Value = 6000
$ws.Range("L44").Value = 7000
$ws.Range("M44").Value = 6500
$ws.Range('N44').Value = '$/docena de matas'
$ws.Range('O44').Value = 'Región de Coquimbo'
$ws.Range("P44").Value = 1083
$ws.Range("Q44").Value = 6
$ws.Range('R44').Value = 'Hortaliza'
